# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value, for each affected sheet.
$updates1 = @{
    3  = 176
    5  = 4976
    9  = 545
    10 = 508
    11 = 1034
    13 = 1378
    14 = 3638
    15 = 409
    17 = 117
    18 = 79
    19 = 2627
    21 = 12
    24 = 176
    25 = 50
    27 = 58
    28 = 266
}

$updates4 = @{
    3  = 176
    6  = 4976
    10 = 545
    11 = 508
    12 = 1034
    14 = 1378
    15 = 3638
    16 = 409
    18 = 117
    19 = 79
    20 = 2627
    22 = 12
    25 = 176
    26 = 50
    28 = 58
    29 = 266
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
